# Corrige a ordenação da coluna "Ano" (C) para o cálculo do índice deflator:
# cada data é deslocada um ano para trás e o respectivo valor (D) é
# recalculado invertendo a base da variação percentual.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $dateText = [string]$cCell.Value2
    $parts = $dateText.Split('/')
    $day = $parts[0]
    $month = $parts[1]
    $year = [int]$parts[2]
    $newYear = $year - 1
    $newDateText = "$day/$month/$newYear"

    $oldValue = [double]$dCell.Value2
    $newValue = (1.0 / (1.0 + $oldValue / 100.0) - 1.0) * 100.0

    # write the new year text while preventing Excel from auto-converting it
    # into a date serial number, then restore the default (unstyled) format
    $cCell.NumberFormat = "@"
    $cCell.Value = $newDateText
    $cCell.Style = "Normal"

    $dCell.Value = $newValue
}
